$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks first (Range.Hyperlinks.Delete clears the whole sheet's collection
# in this COM shim, so one call suffices).
$ws.Range("A1:G30").Hyperlinks.Delete()

# Fully clear old data rows (keep header row 1 untouched).
$ws.Range("A2:G30").Clear()

$data = @{
    2 = @{ 'A'='C2'; 'B'='Stuff'; 'C'=1 }
    3 = @{ 'A'='C2'; 'B'='UX100'; 'C'=5; 'D'='UX100' }
    4 = @{ 'A'='E4'; 'B'='Touchstone Stuff'; 'C'=1; 'G'='https://www.saturdayeveningpost.com/wp-content/uploads/satevepost/2019-12-19-random-stuff-860x573.jpg' }
    5 = @{ 'A'='B1'; 'B'='Mayfield Stuff'; 'C'=1; 'G'='https://media-ecn.s3.amazonaws.com/embedded_image/2016/02/fda.jpg' }
    6 = @{ 'A'='B2'; 'B'='Omni-Tract Stuff'; 'C'=1; 'G'='https://products.integralife.com/ccstore/v1/images/?source=/file/products/Omni-Tract%20Flexible%20Wishbone%20Urologic%20Surgery%20Retractor%20System%20OS%201%20Image.png' }
    7 = @{ 'A'='A1'; 'B'='BNS RF Lesion Generator for Neurosurgery'; 'C'=1; 'D'='RFE2-C'; 'G'='https://www.bnsmed.com/data/watermark/20200924/5f6c31aea1382.jpg' }
    8 = @{ 'A'='A1'; 'B'='Codman Electrosurgical Generator'; 'C'=1; 'D'='901001ESUO'; 'G'='https://products.integralife.com/ccstore/v1/images/?source=/file/products/Codman%20Electrosurgical%20Generator%20OS%201%20Image.jpg' }
    9 = @{ 'A'='A1'; 'B'='Elliquence Surgi-Max Plus'; 'C'=1; 'D'='IEC4-SP'; 'G'='https://www.elliquence.com/wp-content/uploads/2016/01/Surgi-Max-Plus-Device.jpg' }
    10 = @{ 'A'='A2'; 'B'='Integra Duo Headlight & Accessory'; 'C'=1; 'D'='90600'; 'G'='https://www.aamedicalstore.com/SSP Applications/AA Medical SCA/AA Medical/img/Product Images/Integra-Duo-LED-Headlight-Set_01.JPG' }
    11 = @{ 'A'='A2'; 'B'='Lextec Lightsource'; 'C'=1; 'D'='00MLX'; 'G'='https://products.integralife.com/ccstore/v1/images/?source=/file/v6400991064904479991/products/MLX-300-Xenon-Lightsources.jpg' }
    12 = @{ 'A'='A3'; 'B'='BNS 4-Channel RF Lesion Generator'; 'C'=1; 'D'='RFE4-B'; 'G'='https://www.bnsmed.com/data/watermark/20200924/5f6c30bda627b.jpg' }
    13 = @{ 'A'='A3'; 'B'='BNS RF Lesion Generator for Neurosurgery'; 'C'=1; 'D'='RFE2-C'; 'G'='https://www.bnsmed.com/data/watermark/20200924/5f6c31aea1382.jpg' }
    14 = @{ 'A'='C1'; 'B'='Codman Certas Plus'; 'C'=1; 'D'='82-8852'; 'F'='System Failure, Missing Magnet'; 'G'='https://products.integralife.com/ccstore/v1/images/?source=/file/v3841902670343812321/products/ETK_01.png' }
    15 = @{ 'A'='C1'; 'B'='Codman Certas Plus'; 'C'=1; 'D'='82-8852'; 'F'='Unable to power-on'; 'G'='https://products.integralife.com/ccstore/v1/images/?source=/file/v3841902670343812321/products/ETK_01.png' }
    16 = @{ 'A'='C1'; 'B'='Codman Licox PtO2 Monitor'; 'C'=1; 'D'='LCX02'; 'E'='2150601326'; 'F'='Functional'; 'G'='https://products.integralife.com/ccstore/v1/images/?source=/file/v7357354864197611707/collections/licox.jpg' }
    17 = @{ 'A'='C1'; 'B'='Codman Medos Valve Programmer'; 'C'=1; 'D'='82-3126'; 'E'='847'; 'F'='Functional'; 'G'='https://products.integralife.com/ccstore/v1/images/?source=/file/v5137398853523069574/products/823190.jpg' }
    18 = @{ 'A'='C1'; 'B'='Codman Medos Valve Programmer'; 'C'=1; 'D'='82-3126'; 'E'='1173'; 'F'='Functional'; 'G'='https://products.integralife.com/ccstore/v1/images/?source=/file/v5137398853523069574/products/823190.jpg' }
    19 = @{ 'A'='C1'; 'B'='Integra LicocCMP Tissue Oxygen Pressure Monitor'; 'C'=1; 'D'='144733'; 'E'='1629'; 'F'='Missing Power Supply'; 'G'='https://products.integralife.com/ccstore/v1/images/?source=/file/v7357354864197611707/collections/licox.jpg' }
    20 = @{ 'A'='C1'; 'B'='Integra Luxtec Lightsource'; 'C'=1; 'D'='00MLX'; 'E'='16G00MLX7347'; 'F'='GHK Trade-in, Dead motherboard'; 'G'='https://products.integralife.com/ccstore/v1/images/?source=/file/v6400991064904479991/products/MLX-300-Xenon-Lightsources.jpg' }
    21 = @{ 'A'='C1'; 'B'='Integra Luxtec Lightsource'; 'C'=1; 'D'='00MLX'; 'E'='16K00MLX7896'; 'F'='STH Trade-in, Dead motherboard'; 'G'='https://products.integralife.com/ccstore/v1/images/?source=/file/v6400991064904479991/products/MLX-300-Xenon-Lightsources.jpg' }
    22 = @{ 'A'='B3'; 'B'='Stuff'; 'C'=1; 'G'='https://www.saturdayeveningpost.com/wp-content/uploads/satevepost/2019-12-19-random-stuff-860x573.jpg' }
    23 = @{ 'A'='B3'; 'B'='Test'; 'C'=1 }
    24 = @{ 'A'='D2'; 'B'='New Item'; 'C'=1 }
    25 = @{ 'A'='D2'; 'B'='New Item'; 'C'=1 }
}

$cols = @('A','B','C','D','E','F','G')
foreach ($r in ($data.Keys | Sort-Object)) {
    $row = $data[$r]
    foreach ($col in $cols) {
        if ($row.ContainsKey($col)) {
            $ws.Range("$col$r").Value = $row[$col]
        }
    }
}

# Re-create hyperlinks for column G, rows 4..22, in order (matches the target rId sequence).
$gRows = 4..22
foreach ($r in $gRows) {
    $url = $data[$r]['G']
    if ($url) {
        $ws.Hyperlinks.Add($ws.Range("G$r"), $url)
        $ws.Range("G$r").Style = "Hyperlink"
    }
}

Write-Host "done"
